$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")
Write-Host $ws.Name
